$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the added columns
$ws.Range("D1").Value = "freq_2"
$ws.Range("E1").Value = "help_2"

# New data row 10 (session 12)
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 5

# Match the final selection shown in the diff
$ws.Range("E9").Select()
